# Update cryptos list worksheet - mirrors the source data refresh performed
# by the GitHub Actions workflow.
#
# Many "Price" cells contain values that look numeric to Excel's auto-typing
# (e.g. "240.84"), but the source data keeps them as plain text. To reproduce
# that faithfully via COM we briefly mark the cell as Text ("@") before
# assigning the value, then reset the cell style back to Normal so we don't
# leave a stray number format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "36.393.26"
$ws.Range("E2").Value = "  +0.11%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.932.75"
$ws.Range("E3").Value = "  -2.20%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.17%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "240.84"
$ws.Range("E5").Value = "  -1.78%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.605"
$ws.Range("E6").Value = "  -3.60%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - Solana
Set-TextValue $ws.Range("D8") "56.55"
$ws.Range("E8").Value = "  -4.88%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.357"
$ws.Range("E9").Value = "  -4.88%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.02%  "

# Row 11 - TRON
Set-TextValue $ws.Range("D11") "0.103"
$ws.Range("E11").Value = "  -0.72%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.216.26"
$ws.Range("E12").Value = "  -2.20%  "

# Row 13 - Polygon
Set-TextValue $ws.Range("D13") "0.797"
$ws.Range("E13").Value = "  -7.55%  "

# Row 14 - Chainlink
Set-TextValue $ws.Range("D14") "13.32"
$ws.Range("E14").Value = "  -4.50%  "

# Row 15 - Avalanche
Set-TextValue $ws.Range("D15") "20.78"
$ws.Range("E15").Value = "  -11.83%  "

# Row 16 - Polkadot
Set-TextValue $ws.Range("D16") "5.11"
$ws.Range("E16").Value = "  -6.34%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "1.938.21"
$ws.Range("E17").Value = "  -2.01%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "36.290.94"
$ws.Range("E18").Value = "  +0.25%  "

# Row 19 - Litecoin
Set-TextValue $ws.Range("D19") "68.71"
$ws.Range("E19").Value = "  -1.93%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0857"
$ws.Range("E20").Value = "  -1.53%  "

# Row 21 - BitcoinCash
Set-TextValue $ws.Range("D21") "226.64"
$ws.Range("E21").Value = "  -3.26%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "4.93"
$ws.Range("E22").Value = "  -7.08%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.05%  "

# Row 24 - PancakeSwap
Set-TextValue $ws.Range("D24") "2.32"
$ws.Range("E24").Value = "  -11.32%  "

# Row 25 - Toncoin
Set-TextValue $ws.Range("D25") "2.26"
$ws.Range("E25").Value = "  -2.29%  "

# Row 26 - Cosmos
Set-TextValue $ws.Range("D26") "9.22"
$ws.Range("E26").Value = "  -7.91%  "

# Row 27 - Monero
Set-TextValue $ws.Range("D27") "160.63"
$ws.Range("E27").Value = "  -0.91%  "

# Row 28 - Kaspa
Set-TextValue $ws.Range("D28") "0.129"
$ws.Range("E28").Value = "  -3.25%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "19.09"
$ws.Range("E29").Value = "  -3.66%  "

# Row 30 - Stellar
$ws.Range("E30").Value = "  -2.64%  "

# Row 31 - ImmutableX
Set-TextValue $ws.Range("D31") "1.11"
$ws.Range("E31").Value = "  -6.54%  "

# Row 32 - Filecoin
Set-TextValue $ws.Range("D32") "4.53"
$ws.Range("E32").Value = "  -7.67%  "

# Row 33 - Hedera
Set-TextValue $ws.Range("D33") "0.0623"
$ws.Range("E33").Value = "  -0.66%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D34") "4.13"
$ws.Range("E34").Value = "  -6.06%  "

# Row 35 - BinanceUSD
$ws.Range("E35").Value = "  +0.17%  "

# Row 36 - THORChain
Set-TextValue $ws.Range("D36") "6.04"
$ws.Range("E36").Value = "  -2.82%  "

# Row 37 - WEMIXToken
Set-TextValue $ws.Range("D37") "1.78"
$ws.Range("E37").Value = "  -0.31%  "

# Row 38 - LidoDAOToken
Set-TextValue $ws.Range("D38") "2.11"
$ws.Range("E38").Value = "  -6.90%  "

# Row 39 - RenderToken
Set-TextValue $ws.Range("D39") "2.93"
$ws.Range("E39").Value = "  -3.14%  "

# Row 40 - Cronos
Set-TextValue $ws.Range("D40") "0.0967"
$ws.Range("E40").Value = "  +0.35%  "

# Row 41 - HuobiToken
Set-TextValue $ws.Range("D41") "2.87"
$ws.Range("E41").Value = "  -1.12%  "

# Row 42 - now VeChain (was TrustWalletToken)
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D42") "0.0207"
$ws.Range("E42").Value = "  -3.15%  "

# Row 43 - now TrustWalletToken (was VeChain); note: per source data, the
# link column was NOT updated here (it stays as the VeChain URL) - this
# mirrors an upstream data quirk, so C43 is intentionally left untouched.
$ws.Range("B43").Value = "TrustWalletToken"
Set-TextValue $ws.Range("D43") "1.14"
$ws.Range("E43").Value = "  -8.11%  "

# Row 44 - InjectiveProtocol
Set-TextValue $ws.Range("D44") "15.42"
$ws.Range("E44").Value = "  -4.91%  "

# Row 45 - Maker
$ws.Range("D45").Value = "1.327.30"
$ws.Range("E45").Value = "  -3.03%  "

# Row 46 - ARBITRUM
$ws.Range("E46").Value = "  -7.80%  "

# Row 47 - Aave
Set-TextValue $ws.Range("D47") "85.54"
$ws.Range("E47").Value = "  -7.39%  "

# Row 48 - FraxShare
Set-TextValue $ws.Range("D48") "7.02"
$ws.Range("E48").Value = "  -6.27%  "

# Row 49 - MXToken
$ws.Range("E49").Value = "  -0.76%  "

# Row 50 - MultiversX
Set-TextValue $ws.Range("D50") "43.71"
$ws.Range("E50").Value = "  -3.85%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "2.108.83"
$ws.Range("E51").Value = "  -2.25%  "
